function Set-RowValues($sheet, $addr, $count, $vals) {
    $arr = New-Object 'object[,]' 1,$count
    for ($i = 0; $i -lt $count; $i++) { $arr[0,$i] = $vals[$i] }
    $sheet.Range($addr).Value = $arr
}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

Set-RowValues $ws "F366:V366" 17 @(14, 6, -2, 0, -5, -9, 64, 41, 25, 1019, 1016, 1013, 10, 10, 10, 29, 16)
Set-RowValues $ws "X366:Y366" 2 @(0, 3)
$ws.Range("AA366").Value = 209
Set-RowValues $ws "F367:V367" 17 @(17, 12, 5, 3, -1, -3, 87, 45, 26, 1016, 1014, 1013, 10, 9, 8, 32, 16)
Set-RowValues $ws "X367:AA367" 4 @(0, 3, "Rain", 244)
Set-RowValues $ws "F368:V368" 17 @(11, 6, 2, 3, 0, -5, 100, 68, 35, 1016, 1014, 1011, 10, 10, 8, 26, 11)
Set-RowValues $ws "X368:AA368" 4 @(0, 3, "Rain", 295)
Set-RowValues $ws "F369:V369" 17 @(11, 6, 0, -1, -4, -7, 87, 51, 28, 1021, 1017, 1015, 10, 10, 10, 32, 18)
Set-RowValues $ws "X369:Y369" 2 @(0, 2)
$ws.Range("AA369").Value = 297
Set-RowValues $ws "F370:V370" 17 @(17, 7, -2, 0, -3, -4, 86, 56, 27, 1024, 1023, 1021, 10, 10, 10, 19, 8)
Set-RowValues $ws "X370:Y370" 2 @(0, 6)
$ws.Range("AA370").Value = 351
Set-RowValues $ws "F371:Q371" 12 @(18, 9, 1, 1, -3, -7, 75, 43, 18, 1024, 1022, 1020)
Set-RowValues $ws "U371:V371" 2 @(26, 10)
$ws.Range("X371").Value = 0
$ws.Range("AA371").Value = 177
Set-RowValues $ws "F372:V372" 17 @(13, 11, 9, 3, 0, -5, 66, 45, 29, 1022, 1020, 1018, 10, 10, 8, 35, 19)
Set-RowValues $ws "X372:Y372" 2 @(0, 8)
$ws.Range("AA372").Value = 224
Set-RowValues $ws "F373:V373" 17 @(11, 9, 8, 5, 4, 2, 76, 69, 54, 1017, 1013, 1010, 10, 9, 8, 35, 19)
Set-RowValues $ws "X373:AA373" 4 @(0, 7, "Rain-Thunderstorm", 190)
Set-RowValues $ws "F374:V374" 17 @(13, 8, 4, 7, 3, 0, 100, 70, 47, 1014, 1012, 1008, 10, 8, 4, 32, 16)
Set-RowValues $ws "X374:AA374" 4 @(0, 5, "Rain", 288)
Set-RowValues $ws "F375:V375" 17 @(12, 7, 3, 4, -1, -3, 100, 58, 35, 1014, 1012, 1009, 10, 9, 4, 32, 18)
Set-RowValues $ws "X375:AA375" 4 @(0, 4, "Rain", 243)
Set-RowValues $ws "F376:V376" 17 @(11, 6, 1, -1, -4, -6, 75, 47, 30, 1018, 1017, 1014, 10, 10, 10, 32, 16)
Set-RowValues $ws "X376:Y376" 2 @(0, 2)
$ws.Range("AA376").Value = 283
Set-RowValues $ws "F377:V377" 17 @(13, 4, -3, 0, -4, -9, 86, 62, 22, 1018, 1017, 1013, 10, 10, 7, 26, 8)
Set-RowValues $ws "X377:Y377" 2 @(0, 4)
$ws.Range("AA377").Value = 355
Set-RowValues $ws "F378:V378" 17 @(3, 1, 0, 1, -1, -3, 100, 87, 70, 1021, 1018, 1015, 10, 7, 2, 26, 10)
Set-RowValues $ws "X378:AA378" 4 @(0, 7, "Rain-Snow", 324)
Set-RowValues $ws "F379:V379" 17 @(11, 3, -5, -1, -4, -7, 93, 62, 30, 1023, 1021, 1019, 10, 10, 10, 14, 5)
Set-RowValues $ws "X379:Y379" 2 @(0, 1)
$ws.Range("AA379").Value = 298
Set-RowValues $ws "F380:V380" 17 @(14, 6, -2, -2, -4, -8, 81, 49, 23, 1020, 1018, 1016, 10, 10, 10, 29, 11)
Set-RowValues $ws "X380:Y380" 2 @(0, 3)
$ws.Range("AA380").Value = 232
Set-RowValues $ws "F381:V381" 17 @(15, 8, 2, -2, -4, -5, 70, 40, 27, 1018, 1017, 1016, 10, 10, 10, 29, 16)
Set-RowValues $ws "X381:Y381" 2 @(0, 2)
$ws.Range("AA381").Value = 275
Set-RowValues $ws "F382:V382" 17 @(17, 9, 1, -1, -3, -4, 75, 43, 26, 1018, 1017, 1015, 10, 10, 10, 23, 13)
Set-RowValues $ws "X382:Y382" 2 @(0, 3)
$ws.Range("AA382").Value = 272
Set-RowValues $ws "F383:V383" 17 @(18, 9, 1, 0, -3, -6, 81, 45, 19, 1018, 1017, 1016, 10, 10, 10, 29, 10)
Set-RowValues $ws "X383:Y383" 2 @(0, 2)
$ws.Range("AA383").Value = 294
Set-RowValues $ws "F384:V384" 17 @(18, 10, 2, 3, -1, -4, 81, 47, 26, 1020, 1019, 1017, 10, 10, 10, 26, 8)
Set-RowValues $ws "X384:Y384" 2 @(0, 3)
$ws.Range("AA384").Value = 27
Set-RowValues $ws "F385:V385" 17 @(17, 10, 3, 4, 1, -2, 76, 52, 29, 1022, 1020, 1018, 10, 10, 10, 26, 6)
Set-RowValues $ws "X385:Y385" 2 @(0, 3)
$ws.Range("AA385").Value = 313
Set-RowValues $ws "F386:V386" 17 @(17, 9, 2, 4, 2, 1, 93, 60, 34, 1021, 1020, 1018, 10, 10, 10, 29, 8)
Set-RowValues $ws "X386:Y386" 2 @(0, 4)
$ws.Range("AA386").Value = 109
Set-RowValues $ws "F387:V387" 17 @(17, 9, 2, 5, 3, -1, 87, 59, 39, 1021, 1020, 1019, 10, 10, 10, 32, 8)
Set-RowValues $ws "X387:Y387" 2 @(0, 6)
$ws.Range("AA387").Value = 329
Set-RowValues $ws "F388:V388" 17 @(15, 9, 4, 9, 6, 4, 100, 76, 51, 1021, 1019, 1016, 10, 9, 4, 26, 11)
Set-RowValues $ws "X388:AA388" 4 @(0, 6, "Rain", 278)

$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 157
$win.ScrollColumn = 1
$ws.Range("G170").Select()
